$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.218.86"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4835"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2869"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06579"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.890.28"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.72"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07318"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.123"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.99"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6537"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.209.20"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9994"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007747"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.125.52"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.368"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "193.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -7.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.114"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.53"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.05%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.59%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.259"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09072"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.004"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05055"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7133"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.095"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.81%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01779"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.635"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9213"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.043"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "105.74"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.777"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9993"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.373"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1307"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.85"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.930"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05755"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.63"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3806"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.30%  "
